$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.955666
$ws.Range("H2").Value = 2.866998
$ws.Range("I2").Value = 0.1700560254029595
$ws.Range("J2").Value = 0.1700560254029595
$ws.Range("M2").Value = 3.456265333333333
$ws.Range("N2").Value = 10.368796
$ws.Range("O2").Value = 0.009841535807677501
$ws.Range("P2").Value = 0.0098415358076775
$ws.Range("Q2").Value = 3.303035266045333
$ws.Range("R2").Value = 29.727317394408
$ws.Range("S2").Value = 0.001673612463314541
$ws.Range("T2").Value = 0.00167361246331454
$ws.Range("G3").Value = 0.955666
$ws.Range("H3").Value = 2.866998
$ws.Range("I3").Value = 0.1700560254029595
$ws.Range("J3").Value = 0.1700560254029595
$ws.Range("O3").Value = 0.8587907398420774
$ws.Range("P3").Value = 0.8587907398420773
$ws.Range("Q3").Value = 288.2290076756786
$ws.Range("R3").Value = 2594.061069081108
$ws.Range("S3").Value = 0.1460425398704107
$ws.Range("T3").Value = 0.1460425398704107
$ws.Range("G4").Value = 0.955666
$ws.Range("H4").Value = 2.866998
$ws.Range("I4").Value = 0.1700560254029595
$ws.Range("J4").Value = 0.1700560254029595
$ws.Range("O4").Value = 0.1313677243502452
$ws.Range("P4").Value = 0.1313677243502452
$ws.Range("Q4").Value = 44.08988950794468
$ws.Range("R4").Value = 396.809005571502
$ws.Range("S4").Value = 0.02233987306923428
$ws.Range("T4").Value = 0.02233987306923427
$ws.Range("I5").Value = 0.6638424218367511
$ws.Range("J5").Value = 0.663842421836751
$ws.Range("M5").Value = 3.456265333333333
$ws.Range("N5").Value = 10.368796
$ws.Range("O5").Value = 0.009841535807677501
$ws.Range("P5").Value = 0.0098415358076775
$ws.Range("Q5").Value = 12.89395612550622
$ws.Range("R5").Value = 116.045605129556
$ws.Range("S5").Value = 0.006533228965161739
$ws.Range("T5").Value = 0.006533228965161737
$ws.Range("I6").Value = 0.6638424218367511
$ws.Range("J6").Value = 0.663842421836751
$ws.Range("O6").Value = 0.8587907398420774
$ws.Range("P6").Value = 0.8587907398420773
$ws.Range("S6").Value = 0.57010172458774
$ws.Range("T6").Value = 0.5701017245877398
$ws.Range("I7").Value = 0.6638424218367511
$ws.Range("J7").Value = 0.663842421836751
$ws.Range("O7").Value = 0.1313677243502452
$ws.Range("P7").Value = 0.1313677243502452
$ws.Range("S7").Value = 0.08720746828384951
$ws.Range("T7").Value = 0.08720746828384947
$ws.Range("G8").Value = 0.9334429999999999
$ws.Range("I8").Value = 0.1661015527602894
$ws.Range("J8").Value = 0.1661015527602894
$ws.Range("M8").Value = 3.456265333333333
$ws.Range("N8").Value = 10.368796
$ws.Range("O8").Value = 0.009841535807677501
$ws.Range("P8").Value = 0.0098415358076775
$ws.Range("Q8").Value = 3.226226681542666
$ws.Range("R8").Value = 29.03604013388399
$ws.Range("S8").Value = 0.001634694379201222
$ws.Range("T8").Value = 0.001634694379201221
$ws.Range("G9").Value = 0.9334429999999999
$ws.Range("I9").Value = 0.1661015527602894
$ws.Range("J9").Value = 0.1661015527602894
$ws.Range("O9").Value = 0.8587907398420774
$ws.Range("P9").Value = 0.8587907398420773
$ws.Range("S9").Value = 0.1426464753839268
$ws.Range("T9").Value = 0.1426464753839267
$ws.Range("G10").Value = 0.9334429999999999
$ws.Range("I10").Value = 0.1661015527602894
$ws.Range("J10").Value = 0.1661015527602894
$ws.Range("O10").Value = 0.1313677243502452
$ws.Range("P10").Value = 0.1313677243502452
$ws.Range("Q10").Value = 43.06462585460233
$ws.Range("R10").Value = 387.581632691421
$ws.Range("S10").Value = 0.0218203829971614
$ws.Range("T10").Value = 0.0218203829971614
